{"js": "// Retitle the document and switch the \"My question / hypothesis\" body\n// paragraph from present to past tense, per the authoring diff.\n\nconst body = context.document.body;\n\n// --- 1. Title: \"Final Report\" -> \"The Negative Effect of Air Quality on Birds\"\nconst titleParagraph = body.paragraphs.getFirst();\ntitleParagraph.load(\"style\");\nawait context.sync();\n\nif (titleParagraph.style === \"Title\") {\n  titleParagraph.insertText(\n    \"The Negative Effect of Air Quality on Birds\",\n    \"Replace\"\n  );\n} else {\n  // Fall back to a scoped search in case paragraph ordering ever changes.\n  const titleHits = body.search(\"Final Report\", { matchCase: true });\n  titleHits.load(\"items\");\n  await context.sync();\n  titleHits.items[0].insertText(\n    \"The Negative Effect of Air Quality on Birds\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// --- 2. Body paragraph tense fixes (each phrase is unique in the document).\nconst replacements = [\n  [\"My question is\", \"My question was\"],\n  [\"My hypothesis is\", \"My hypothesis was\"],\n  [\"Another question I have is\", \"Another question I had was\"],\n  [\"I believe that there is\", \"I believed that there was\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const hits = body.search(findText, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(replaceText, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Retitle the document and switch the \"My question / hypothesis\" body\n# paragraph from present to past tense, per the authoring diff.\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceOne = 2\n    $find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2) | Out-Null\n}\n\n# --- 1. Title: \"Final Report\" -> \"The Negative Effect of Air Quality on Birds\"\nReplace-FirstMatch \"Final Report\" \"The Negative Effect of Air Quality on Birds\"\n\n# --- 2. Body paragraph tense fixes (each phrase is unique in the document).\nReplace-FirstMatch \"My question is\" \"My question was\"\nReplace-FirstMatch \"My hypothesis is\" \"My hypothesis was\"\nReplace-FirstMatch \"Another question I have is\" \"Another question I had was\"\nReplace-FirstMatch \"I believe that there is\" \"I believed that there was\"\n"}
